# Regenerate merged AHB files
#
# The sheet "AHB-Diff" has a header row (row 1, A1:U59 data range) whose
# column captions were produced by an older merge ("_old" / "_new" suffixes).
# This regenerates the header captions to reflect the actual compared
# versions (FV2310 / FV2404), wraps the range in a proper Excel Table
# ("Table1") so the column names are addressable, and freezes the header
# row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header captions -----------------------------------
# Columns A:J were "<Name>_old"  -> "<Name>_FV2310"
# Column  K        stays "diff"
# Columns L:U were "<Name>_new"  -> "<Name>_FV2404"
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2310"
}
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2404"
}

# --- 2. Turn the used range into an Excel Table ------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U59"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
